$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1107347873458729
$ws.Range("C2").Value = 0.002963102129239998

$ws.Range("B3").Value = 0.1076506458811195
$ws.Range("C3").Value = 0.001400656513985915

$ws.Range("B4").Value = 0.08703852422468288
$ws.Range("C4").Value = 0.002767699295926451

$ws.Range("B5").Value = 0.09880964859350067
$ws.Range("C5").Value = 0.003304587644529865

$ws.Range("B6").Value = 0.05622964218543108
$ws.Range("C6").Value = 0.001546233530488934

$ws.Range("B7").Value = 0.02994570698494851
$ws.Range("C7").Value = 0.002914276493481508

$ws.Range("B8").Value = 0.01071931590666068
$ws.Range("C8").Value = 0.0008408392436984117
